$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TC-25 / TC-26 added in TestPlan tab.
# tc024 already holds Project Name / release / Status columns with the
# correct formatting (style s="4" on columns A/B) - duplicate it twice so
# the new tc025 / tc026 sheets inherit that formatting, column widths, etc,
# then patch the per-sheet differences (extra columns, header/data text,
# selection, tab order).
# ---------------------------------------------------------------------------

$tc024 = $wb.Worksheets.Item("tc024")

# tc024: fix the typo in the release-update-notification string and move the
# selection from C2 to B2 (tabSelected is dropped from tc024 automatically
# once another sheet becomes active later in this script).
$tc024.Range("B2").Value = " Release update notification 09-01-2026"
$tc024.Range("B2").Select()

# --- tc025 ------------------------------------------------------------------
$tc024.Copy([System.Reflection.Missing]::Value, $tc024)
$tc025 = $wb.Worksheets.Item("tc024 (2)")
$tc025.Name = "tc025"

$tc025.Range("C1").Value = "Tesplaname"
$tc025.Range("D1").Value = "des"
$tc025.Range("C2").Value = "Cyle update notification 09-01-2026"
$tc025.Range("D2").Value = "planned"

$tc025.Columns.Item(1).ColumnWidth = 21.98307291666667
$tc025.Columns.Item(2).ColumnWidth = 33.79947916666666
$tc025.Columns.Item(3).ColumnWidth = 30.34635416666667

$tc025.Range("A1:D2").Select()

# --- tc026 ------------------------------------------------------------------
$tc025.Copy([System.Reflection.Missing]::Value, $tc025)
$tc026 = $wb.Worksheets.Item("tc025 (2)")
$tc026.Name = "tc026"

$tc026.Range("C1").Value = "Cyclename"
$tc026.Range("D1").Value = "Suitename"
$tc026.Range("E1").Value = "Desc"
$tc026.Range("C2").Value = "Cyle update notification 09-01-2026"
$tc026.Range("D2").Value = "Suite update notification 09-01-2026"
$tc026.Range("E2").Value = "Planned"

$tc026.Columns.Item(2).ColumnWidth = 37.70963541666666
$tc026.Columns.Item(3).ColumnWidth = 30.79947916666667
$tc026.Columns.Item(4).ColumnWidth = 37.52994791666666

$tc026.Range("D6").Select()
